$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand-new row just above the totals row (old row 62), which
#    pushes the totals row to 63 and the footer row to 64. This gives us a
#    blank row 62 that we turn into a duplicate of the (old) last item row.
# ---------------------------------------------------------------------------
$ws.Rows("62:62").Insert()

# Copy the formatting (styles, fonts, borders, number formats) from row 61
# into the freshly inserted row 62, then fix up the row height and merges
# which a plain "insert" does not carry over in this engine.
$ws.Range("A61:Q61").Copy()
$ws.Range("A62:Q62").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Rows("62:62").RowHeight = $ws.Rows("61:61").RowHeight

$ws.Range("A62:B62").Merge()
$ws.Range("C62:G62").Merge()
$ws.Range("H62:K62").Merge()
$ws.Range("L62:M62").Merge()
$ws.Range("N62:O62").Merge()

# ---------------------------------------------------------------------------
# 2) The new row 62 becomes a copy of what used to be the last item
#    (row 61, "معجون سيجنال 50 مل") - i.e. every item row from 50-61 is about
#    to shift down by one logical position, with row 62 taking on the values
#    that used to live in row 61.
# ---------------------------------------------------------------------------
function Set-TextValue($rng, $text) {
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.NumberFormat = $fmt
}

$ws.Range("A62").Value2 = 56
Set-TextValue $ws.Range("C62") "معجون سيجنال 50 مل"
Set-TextValue $ws.Range("H62") "3:0"
Set-TextValue $ws.Range("L62") "0"
Set-TextValue $ws.Range("N62") "35.00"
Set-TextValue $ws.Range("P62") "35.0000"
Set-TextValue $ws.Range("Q62") "1:0"

# ---------------------------------------------------------------------------
# 3) Shift the remaining item rows (50-60 -> 51-61) down one position by
#    copying each row's text fields into the row below it, working from the
#    bottom up so nothing is overwritten before it is copied.
# ---------------------------------------------------------------------------
$items = @(
  @{H="2:0";  N="45.00"; P="45.0000"; Q="1:0"},  # was row 50 -> now 51
  @{H="0:0";  N="4.00";  P="8.0000";  Q="2:0"},   # was row 51 -> now 52
  @{H="0:0";  N="2.00";  P="8.0000";  Q="4:0"},   # was row 52 -> now 53
  @{H="0:0";  N="3.00";  P="9.0000";  Q="3:0"},   # was row 53 -> now 54
  @{H="0:0";  N="4.00";  P="4.0000";  Q="1:0"},   # was row 54 -> now 55
  @{H="5:0";  N="20.00"; P="20.0000"; Q="1:0"},   # was row 55 -> now 56
  @{H="2:0";  N="12.00"; P="12.0000"; Q="1:0"},   # was row 56 -> now 57
  @{H="0:0";  N="65.00"; P="65.0000"; Q="1:0"},   # was row 57 -> now 58
  @{H="0:0";  N="15.00"; P="30.0000"; Q="2:0"},   # was row 58 -> now 59
  @{H="1:0";  N="24.00"; P="24.0000"; Q="1:0"},   # was row 59 -> now 60
  @{H="28:0"; N="24.00"; P="24.0000"; Q="1:0"}    # was row 60 -> now 61
)
$names = @(
  "حنه جلوري سوده 1 كيس",
  "سرنجات 10 سم",
  "سرنجات 3 سم",
  "سرنجات 5 سم",
  "سرنجه دواء",
  "سلاكه ازن خشب كبيره",
  "شاش فازلين 10*20 سم ",
  "شفاط ثدي فلاي بيبي",
  "كالونا ",
  "محلول رينجر",
  "محلول ملح"
)

for ($i = 10; $i -ge 0; $i--) {
    $destRow = 51 + $i
    Set-TextValue $ws.Range("C$destRow") $names[$i]
    Set-TextValue $ws.Range("H$destRow") $items[$i].H
    Set-TextValue $ws.Range("N$destRow") $items[$i].N
    Set-TextValue $ws.Range("P$destRow") $items[$i].P
    Set-TextValue $ws.Range("Q$destRow") $items[$i].Q
}

# ---------------------------------------------------------------------------
# 4) Row 50 becomes the new item "حبايه".
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("C50") "حبايه"
Set-TextValue $ws.Range("H50") "0:0"
Set-TextValue $ws.Range("N50") "3.00"
Set-TextValue $ws.Range("P50") "12.0000"
Set-TextValue $ws.Range("Q50") "4:0"

# ---------------------------------------------------------------------------
# 5) "سرنجات 3 سم" (now row 53) independently changed its selling price and
#    transaction count.
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("P53") "18.0000"
Set-TextValue $ws.Range("Q53") "9:0"

# ---------------------------------------------------------------------------
# 6) Update the grand total (now row 63) to reflect the new/changed rows.
# ---------------------------------------------------------------------------
$ws.Range("P63").Value2 = 3183.09
$ws.Rows("63:63").RowHeight = 24.75

# ---------------------------------------------------------------------------
# 7) Update the generated-on timestamp in the footer (now row 64).
# ---------------------------------------------------------------------------
$ws.Range("A64").Value2 = "Tuesday, 29 July, 2025 5:42 PM"
